$d = $word.ActiveDocument

$pairs = @(
    @("961×7=", "938×9="),
    @("912×2=", "744×3="),
    @("533×3=", "230×4="),
    @("881×3=", "113×3="),
    @("905×9=", "732×3="),
    @("844×2=", "678×4="),
    @("680×5=", "432×2="),
    @("265×3=", "984×7="),
    @("136×4=", "288×2="),
    @("251×6=", "783×8="),
    @("861×2=", "739×4="),
    @("901×8=", "635×3="),
    @("675×2=", "281×4="),
    @("509×8=", "646×7="),
    @("301×5=", "825×4="),
    @("999×3=", "777×2="),
    @("622×8=", "235×2="),
    @("828×4=", "629×3="),
    @("197×8=", "918×7="),
    @("534×8=", "475×5="),
    @("963×4=", "376×9="),
    @("419×9=", "242×2="),
    @("109×9=", "630×8="),
    @("622×6=", "757×8="),
    @("330×7=", "297×2=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
